# Scheduled-runner price/profit refresh: update cached market-board figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4684.5
$ws.Range("I43").Value = 4159.6665
$ws.Range("J43").Value = 4999.4
$ws.Range("K43").Value = 4159.6665
$ws.Range("L43").Value = 4999.4
$ws.Range("M43").Value = -4090.6665
$ws.Range("N43").Value = -5137.4

$ws.Range("H70").Value = 1486.6666
$ws.Range("I70").Value = 1192
$ws.Range("K70").Value = 3576
$ws.Range("M70").Value = -3306

$ws.Range("H73").Value = 1486.6666
$ws.Range("I73").Value = 1192
$ws.Range("K73").Value = 3576
$ws.Range("M73").Value = -2640

$ws.Range("H100").Value = 8374.286
$ws.Range("I100").Value = 8374.286
$ws.Range("K100").Value = 8374.286
$ws.Range("M100").Value = -7833.286

$ws.Range("H116").Value = 2938.8572
$ws.Range("I116").Value = 2699.3333
$ws.Range("J116").Value = 3118.5
$ws.Range("K116").Value = 2699.3333
$ws.Range("L116").Value = 3118.5
$ws.Range("M116").Value = 742.6667000000002
$ws.Range("N116").Value = -10002.5

$ws.Range("H118").Value = 18518848
$ws.Range("I118").Value = 20833674
$ws.Range("J118").Value = 250
$ws.Range("K118").Value = 62501022
$ws.Range("L118").Value = 750
$ws.Range("M118").Value = -62499365
$ws.Range("N118").Value = -4064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2846.851
$ws.Range("I32").Value = 2821.8044
$ws.Range("K32").Value = 2821.8044
$ws.Range("M32").Value = -2534.8044

$ws.Range("H34").Value = 25252
$ws.Range("I34").Value = 25252
$ws.Range("K34").Value = 25252
$ws.Range("M34").Value = -24981

$ws.Range("H44").Value = 41999
$ws.Range("I44").Value = 41999
$ws.Range("K44").Value = 41999
$ws.Range("M44").Value = -41511

$ws.Range("H55").Value = 24999.666
$ws.Range("I55").Value = 14999.5
$ws.Range("K55").Value = 14999.5
$ws.Range("M55").Value = -14684.5

$ws.Range("H61").Value = 1883.6666
$ws.Range("I61").Value = 1642.5
$ws.Range("K61").Value = 1642.5
$ws.Range("M61").Value = -1430.5

$ws.Range("H74").Value = 2911.2856
$ws.Range("I74").Value = 2791.5
$ws.Range("K74").Value = 2791.5
$ws.Range("M74").Value = -1917.5

$ws.Range("H77").Value = 2911.2856
$ws.Range("I77").Value = 2791.5
$ws.Range("K77").Value = 13957.5
$ws.Range("M77").Value = -9589.5

$ws.Range("H132").Value = 2245.1875
$ws.Range("J132").Value = 3632.6667
$ws.Range("L132").Value = 10898.0001
$ws.Range("N132").Value = -15958.0001

$ws.Range("H136").Value = 1883.6666
$ws.Range("I136").Value = 1642.5
$ws.Range("K136").Value = 4927.5
$ws.Range("M136").Value = -2377.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 341.33334
$ws.Range("I22").Value = 324.5
$ws.Range("K22").Value = 324.5
$ws.Range("M22").Value = -151.5

$ws.Range("H86").Value = 7125.1875
$ws.Range("I86").Value = 2456.111
$ws.Range("K86").Value = 2456.111
$ws.Range("M86").Value = -1333.111

$ws.Range("H89").Value = 7125.1875
$ws.Range("I89").Value = 2456.111
$ws.Range("K89").Value = 12280.555
$ws.Range("M89").Value = -6664.555

$ws.Range("H134").Value = 4999.5
$ws.Range("I134").Value = 4999.5
$ws.Range("K134").Value = 14998.5
$ws.Range("M134").Value = -12463.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 729.9167
$ws.Range("I16").Value = 650.9
$ws.Range("K16").Value = 650.9
$ws.Range("M16").Value = -363.9

$ws.Range("H50").Value = 27250
$ws.Range("J50").Value = 27250
$ws.Range("L50").Value = 27250
$ws.Range("N50").Value = -28500

$ws.Range("H113").Value = 729.9167
$ws.Range("I113").Value = 650.9
$ws.Range("K113").Value = 650.9
$ws.Range("M113").Value = 1519.1

$ws.Range("H132").Value = 1851.8889
$ws.Range("I132").Value = 1439.6
$ws.Range("J132").Value = 2367.25
$ws.Range("K132").Value = 4318.799999999999
$ws.Range("L132").Value = 7101.75
$ws.Range("M132").Value = -1788.799999999999
$ws.Range("N132").Value = -12161.75

$ws.Range("H134").Value = 4681.107
$ws.Range("I134").Value = 3458.318
$ws.Range("K134").Value = 10374.954
$ws.Range("M134").Value = -7839.954000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 3894.3333
$ws.Range("I21").Value = 119
$ws.Range("J21").Value = 4649.4
$ws.Range("K21").Value = 357
$ws.Range("L21").Value = 13948.2
$ws.Range("M21").Value = -184
$ws.Range("N21").Value = -14294.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 15000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 15000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 15000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -15332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 918
$ws.Range("I22").Value = 1037.4286
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1037.4286
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -742.4286
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 918
$ws.Range("I27").Value = 1037.4286
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 1037.4286
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -930.4286
$ws.Range("N27").Value = -714

$ws.Range("H40").Value = 5396.5
$ws.Range("I40").Value = 5337
$ws.Range("K40").Value = 5337
$ws.Range("M40").Value = -5201

$ws.Range("H46").Value = 3329.6
$ws.Range("I46").Value = 809.4
$ws.Range("J46").Value = 4169.6665
$ws.Range("K46").Value = 809.4
$ws.Range("L46").Value = 4169.6665
$ws.Range("M46").Value = -621.4
$ws.Range("N46").Value = -4545.6665

$ws.Range("H61").Value = 2637.1428
$ws.Range("I61").Value = 2672.6
$ws.Range("K61").Value = 2672.6
$ws.Range("M61").Value = -2470.6

$ws.Range("H113").Value = 2637.1428
$ws.Range("I113").Value = 2672.6
$ws.Range("K113").Value = 2672.6
$ws.Range("M113").Value = -502.5999999999999

$ws.Range("H132").Value = 4010.3914
$ws.Range("I132").Value = 2964.1904
$ws.Range("K132").Value = 8892.5712
$ws.Range("M132").Value = -6362.5712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2759.6191
$ws.Range("I132").Value = 1812.0625
$ws.Range("K132").Value = 5436.1875
$ws.Range("M132").Value = -2906.1875
